$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("grandes regiões e unidades da federação") was a spurious header
# row with no data. Deleting it shifts every subsequent row (7..37) up by
# one, so "norte" (old row 7, with its B:G data) becomes the new row 6,
# "rondônia" becomes row 7, ..., "goiás" (old row 37) becomes row 36.
$ws.Rows("6").Delete()
